$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Set font color for D2 (first cell), creating a new consolidated style
$ws.Range("D2").Font.Color = 0

# 2. Copy that resulting style onto D1, D3, D4, D5 so whole column matches
$ws.Range("D2").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)

# 3. Add color column
$ws.Range("E1:E5").Copy()
$ws.Range("F1:F5").PasteSpecial(-4122)
$ws.Range("F1").Value = "color"
$ws.Range("F2").Value = "['22/255', '39/255', '136/255']"
$ws.Range("F3").Value = "['232/255', '198/255', '28/255']"
$ws.Range("F4").Value = "['72/255', '193/255', '78/255']"
$ws.Range("F5").Value = "['3/255', '97/255', '30/255']"

# 4. Row heights
$ws.Rows("1:5").RowHeight = 19.5
